$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RYB")

# Insert a new column before column D (shifts existing D:K data to E:L)
$ws.Columns("D:D").Insert()

# Copy formatting (number formats, font, etc.) from the now-shifted column E
# into the newly inserted (blank) column D so every row keeps its correct style
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new column D with the latest quarter's figures
$ws.Range("D7").Value = 43373
$ws.Range("D8").Value = 35300
$ws.Range("D9").Value = 34000
$ws.Range("D10").Value = 1300
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 40500
$ws.Range("D18").Value = -5200
$ws.Range("D20").Value = 500
$ws.Range("D21").Value = "NA"
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = -4700
$ws.Range("D24").Value = -800
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = -3800
$ws.Range("D27").Value = -2600
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -500
$ws.Range("D33").Value = -2600
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -2600
$ws.Range("D38").Value = 43373
$ws.Range("D41").Value = 137600
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 1200
$ws.Range("D44").Value = 5000
$ws.Range("D45").Value = 14100
$ws.Range("D46").Value = 157900
$ws.Range("D47").Value = 1800
$ws.Range("D48").Value = 46000
$ws.Range("D49").Value = 31100
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 20600
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 257500
$ws.Range("D57").Value = 0
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 125100
$ws.Range("D60").Value = 125100
$ws.Range("D61").Value = 0
$ws.Range("D62").Value = 20300
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 151000
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -27600
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 106500
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43373
$ws.Range("D81").Value = -2600
$ws.Range("D83").Value = 0
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 0
$ws.Range("D91").Value = 0
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = 0
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 0
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 0
$wb.Save()
